# Weekly fruit/vegetable price update: insert a new daily record row.
#
# The dataset keeps one row per market/date reading. A new reading (row 60,
# dated 44784) is inserted above the former row 60, pushing every
# subsequent row (old 60..88) down by one (new 61..89).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 60 (old rows 60-88 shift down to 61-89,
# carrying the D-column date style along with them).
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with the new reading.
$ws.Range("A60").Value = 11
$ws.Range("B60").Value = "Vega Monumental Concepción"
$ws.Range("C60").Value = "Bíobío"
$ws.Range("D60").Value = 44784
$ws.Range("E60").Value = 8
$ws.Range("F60").Value = 100112012
$ws.Range("G60").Value = "Espinaca"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 40
$ws.Range("K60").Value = 6500
$ws.Range("L60").Value = 7000
$ws.Range("M60").Value = 6750
$ws.Range("N60").Value = "`$/cuna 10 kilos"
$ws.Range("O60").Value = "Región Metropolitana"
$ws.Range("P60").Value = 675
$ws.Range("Q60").Value = 10
$ws.Range("R60").Value = "Hortaliza"
